$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell address -> new value pairs reflecting the latest crypto price/volume snapshot.
$updates = @(
    @{ Cell = 'D2'; Value = '28.202.64' }
    @{ Cell = 'E2'; Value = '  +0.85%  ' }
    @{ Cell = 'D3'; Value = '1.796.32' }
    @{ Cell = 'E3'; Value = '  +2.64%  ' }
    @{ Cell = 'E4'; Value = '  +0.27%  ' }
    @{ Cell = 'D5'; Value = '''331.97' }
    @{ Cell = 'E5'; Value = '  -0.19%  ' }
    @{ Cell = 'E6'; Value = '  +0.47%  ' }
    @{ Cell = 'D7'; Value = '''0.4530' }
    @{ Cell = 'E7'; Value = '  +16.94%  ' }
    @{ Cell = 'D8'; Value = '''0.3721' }
    @{ Cell = 'E8'; Value = '  +9.96%  ' }
    @{ Cell = 'D9'; Value = '''44.96' }
    @{ Cell = 'E9'; Value = '  -1.08%  ' }
    @{ Cell = 'D10'; Value = '''1.142' }
    @{ Cell = 'E10'; Value = '  +2.91%  ' }
    @{ Cell = 'D11'; Value = '''0.07560' }
    @{ Cell = 'E11'; Value = '  +5.40%  ' }
    @{ Cell = 'E12'; Value = '  +0.59%  ' }
    @{ Cell = 'D13'; Value = '''22.34' }
    @{ Cell = 'E13'; Value = '  +0.81%  ' }
    @{ Cell = 'D14'; Value = '''6.301' }
    @{ Cell = 'E14'; Value = '  +3.05%  ' }
    @{ Cell = 'D15'; Value = '''7.501' }
    @{ Cell = 'E15'; Value = '  +7.14%  ' }
    @{ Cell = 'D16'; Value = '1.789.06' }
    @{ Cell = 'E16'; Value = '  +2.53%  ' }
    @{ Cell = 'D17'; Value = '''0.00001090' }
    @{ Cell = 'E17'; Value = '  +3.63%  ' }
    @{ Cell = 'D18'; Value = '''0.06746' }
    @{ Cell = 'E18'; Value = '  +2.45%  ' }
    @{ Cell = 'D19'; Value = '''80.74' }
    @{ Cell = 'E19'; Value = '  +0.61%  ' }
    @{ Cell = 'D20'; Value = '''1.000' }
    @{ Cell = 'E20'; Value = '  +0.40%  ' }
    @{ Cell = 'D21'; Value = '''17.49' }
    @{ Cell = 'E21'; Value = '  +3.59%  ' }
    @{ Cell = 'D22'; Value = '''6.355' }
    @{ Cell = 'E22'; Value = '  +2.95%  ' }
    @{ Cell = 'D23'; Value = '28.196.47' }
    @{ Cell = 'E23'; Value = '  +0.80%  ' }
    @{ Cell = 'D24'; Value = '''11.81' }
    @{ Cell = 'E24'; Value = '  +2.47%  ' }
    @{ Cell = 'D25'; Value = '''2.413' }
    @{ Cell = 'E25'; Value = '  +1.35%  ' }
    @{ Cell = 'D26'; Value = '''20.50' }
    @{ Cell = 'E26'; Value = '  +3.49%  ' }
    @{ Cell = 'B27'; Value = 'LidoDAOToken' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D27'; Value = '''2.367' }
    @{ Cell = 'E27'; Value = '  +2.64%  ' }
    @{ Cell = 'B28'; Value = 'Monero' }
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D28'; Value = '''152.13' }
    @{ Cell = 'E28'; Value = '  -0.73%  ' }
    @{ Cell = 'D29'; Value = '1.993.45' }
    @{ Cell = 'E29'; Value = '  +2.37%  ' }
    @{ Cell = 'D30'; Value = '''132.36' }
    @{ Cell = 'E30'; Value = '  +3.15%  ' }
    @{ Cell = 'D31'; Value = '''1.234' }
    @{ Cell = 'E31'; Value = '  -2.02%  ' }
    @{ Cell = 'D32'; Value = '''4.039' }
    @{ Cell = 'E32'; Value = '  -0.47%  ' }
    @{ Cell = 'D33'; Value = '''0.09434' }
    @{ Cell = 'E33'; Value = '  +9.25%  ' }
    @{ Cell = 'D34'; Value = '''5.800' }
    @{ Cell = 'E34'; Value = '  +0.12%  ' }
    @{ Cell = 'D35'; Value = '''0.2374' }
    @{ Cell = 'E35'; Value = '  +13.47%  ' }
    @{ Cell = 'D36'; Value = '''12.11' }
    @{ Cell = 'E36'; Value = '  +1.01%  ' }
    @{ Cell = 'B37'; Value = 'Hedera' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' }
    @{ Cell = 'D37'; Value = '''0.06322' }
    @{ Cell = 'E37'; Value = '  +3.31%  ' }
    @{ Cell = 'B38'; Value = 'VeChain' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D38'; Value = '''0.02338' }
    @{ Cell = 'E38'; Value = '  +2.90%  ' }
    @{ Cell = 'D39'; Value = '''5.204' }
    @{ Cell = 'E39'; Value = '  +1.79%  ' }
    @{ Cell = 'D40'; Value = '''0.6584' }
    @{ Cell = 'E40'; Value = '  +2.23%  ' }
    @{ Cell = 'D41'; Value = '''8.351' }
    @{ Cell = 'E41'; Value = '  +6.70%  ' }
    @{ Cell = 'D42'; Value = '''1.479' }
    @{ Cell = 'E42'; Value = '  -1.01%  ' }
    @{ Cell = 'D43'; Value = '''1.204' }
    @{ Cell = 'E43'; Value = '  +0.52%  ' }
    @{ Cell = 'B44'; Value = 'EnergySwap' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D44'; Value = '''14.20' }
    @{ Cell = 'E44'; Value = '  +4.18%  ' }
    @{ Cell = 'B45'; Value = 'Frax' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' }
    @{ Cell = 'D45'; Value = '''1.000' }
    @{ Cell = 'E45'; Value = '  +0.35%  ' }
    @{ Cell = 'B46'; Value = 'Decentraland' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D46'; Value = '''0.6108' }
    @{ Cell = 'E46'; Value = '  +2.70%  ' }
    @{ Cell = 'B47'; Value = 'PancakeSwap' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D47'; Value = '''3.809' }
    @{ Cell = 'E47'; Value = '  +0.10%  ' }
    @{ Cell = 'D48'; Value = '''130.05' }
    @{ Cell = 'E48'; Value = '  +3.33%  ' }
    @{ Cell = 'D49'; Value = '''2.032' }
    @{ Cell = 'E49'; Value = '  +3.17%  ' }
    @{ Cell = 'D50'; Value = '''0.07125' }
    @{ Cell = 'E50'; Value = '  +2.04%  ' }
    @{ Cell = 'D51'; Value = '''1.159' }
    @{ Cell = 'E51'; Value = '  +0.84%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

